# Apply updated cryptocurrency price/volume data to worksheet.
# Cells D2:D51 and E2:E51 are plain text cells; a leading apostrophe
# forces Excel to store the value as text instead of auto-converting
# it to a number/date, matching the original inlineStr cell type.
# Style is then reset to "Normal" so no extra quote-prefix styling is
# left behind on the cell (the column has no special number format).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q = [string][char]39

$ws.Range('D2').Value = $q + '23.191.33'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = $q + '  -2.92%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = $q + '1.610.45'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = $q + '  -2.36%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = $q + '0.9989'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = $q + '  -0.26%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = $q + '0.9995'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = $q + '  -0.12%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = $q + '301.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = $q + '  -2.22%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = $q + '0.3785'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = $q + '  -2.72%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = $q + '0.3660'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = $q + '  -4.40%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = $q + '49.56'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = $q + '  -4.19%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = $q + '1.000'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = $q + '  -0.05%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = $q + '1.268'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = $q + '  -6.11%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = $q + '0.08094'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = $q + '  -3.81%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = $q + '23.06'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = $q + '  -3.37%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = $q + '6.614'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = $q + '  -6.54%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = $q + '7.449'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = $q + '  -6.17%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = $q + '0.00001259'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = $q + '  -4.32%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = $q + '1.611.63'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = $q + '  -2.02%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = $q + '91.46'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = $q + '  -3.36%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = $q + '0.06800'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = $q + '  -2.45%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = $q + '18.41'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = $q + '  -6.61%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = $q + '  -5.16%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D23').Value = $q + '13.07'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = $q + '  -4.65%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = $q + '23.182.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = $q + '  -2.93%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = $q + '2.353'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = $q + '  -4.00%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = $q + '2.862'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = $q + '  -2.96%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = $q + '21.10'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Value = $q + '150.66'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = $q + '  -0.55%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = $q + '5.286'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = $q + '  -2.46%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = $q + '132.97'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = $q + '  -4.20%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = $q + '2.416'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = $q + '  -4.03%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = $q + '6.872'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = $q + '  -12.49%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = $q + '1.790.49'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = $q + '  -1.93%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = $q + '0.9708'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = $q + '  -7.27%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = $q + '0.07715'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = $q + '  -4.01%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = $q + '0.02769'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = $q + '  -6.53%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = $q + '0.2563'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = $q + '  -4.48%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = $q + '6.249'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = $q + '  -6.39%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = $q + '10.18'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = $q + '  -7.24%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = $q + '0.08912'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = $q + '  -2.04%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = $q + '1.392'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = $q + '  -2.16%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = $q + '0.7201'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = $q + '  -5.21%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = $q + '12.85'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = $q + '  -4.46%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = $q + '15.74'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = $q + '  -3.16%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = $q + '0.6693'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = $q + '  -4.28%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = $q + '2.309'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = $q + '  -6.28%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = $q + '0.9991'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = $q + '  -0.13%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = $q + '3.981'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = $q + '  -2.26%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = $q + '0.08018'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = $q + '  -3.20%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = $q + '131.05'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = $q + '  -2.52%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = $q + '  -2.91%  '
$ws.Range('E51').Style = 'Normal'
